# Update the "updated_at" timestamps for the rows whose value was
# 2025-12-25 23:25:54 / 2025-12-25 23:25:57 to the new timestamps
# 2025-12-26 00:03:23 / 2025-12-26 00:03:26 respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue1 = "2025-12-25 23:25:54"
$newValue1 = "2025-12-26 00:03:23"
$oldValue2 = "2025-12-25 23:25:57"
$newValue2 = "2025-12-26 00:03:26"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cellValue = $cell.Value()
        if ($cellValue -eq $oldValue1) {
            $cell.Value = $newValue1
        } elseif ($cellValue -eq $oldValue2) {
            $cell.Value = $newValue2
        }
    }
}
